$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2510
$ws.Range("E2").Value = 142
$ws.Range("F2").Value = 142
$ws.Range("G2").Value = 229
$ws.Range("H2").Value = 185
$ws.Range("I2").Value = 185
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 4815
$ws.Range("L2").Value = 286
$ws.Range("M2").Value = 4529
$ws.Range("N2").Value = 4529
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 100
$ws.Range("Q2").Value = 353
$ws.Range("R2").Value = -302
$ws.Range("S2").Value = -44
$ws.Range("T2").Value = 47
$ws.Range("U2").Value = 306
$ws.Range("V2").Value = 1
$ws.Range("W2").Value = 5.68
$ws.Range("X2").Value = 7.38
$ws.Range("Y2").Value = 4.16
$ws.Range("Z2").Value = 3.9
$ws.Range("AA2").Value = 6.31
$ws.Range("AB2").Value = 4167.6
$ws.Range("AC2").Value = 925
$ws.Range("AD2").Value = 12.8
$ws.Range("AE2").Value = 22646
$ws.Range("AF2").Value = 0.52
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 1.69
$ws.Range("AI2").Value = 21.61
$ws.Range("AJ2").Value = 20000000

# Row 3
$ws.Range("D3").Value = 2213
$ws.Range("E3").Value = 117
$ws.Range("F3").Value = 117
$ws.Range("G3").Value = 165
$ws.Range("H3").Value = 138
$ws.Range("I3").Value = 138
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = 4868
$ws.Range("L3").Value = 235
$ws.Range("M3").Value = 4633
$ws.Range("N3").Value = 4633
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = 100
$ws.Range("Q3").Value = 309
$ws.Range("R3").Value = -444
$ws.Range("S3").Value = -39
$ws.Range("T3").Value = 47
$ws.Range("U3").Value = 262
$ws.Range("V3").Value = 3
$ws.Range("W3").Value = 5.28
$ws.Range("X3").Value = 6.25
$ws.Range("Y3").Value = 3.02
$ws.Range("Z3").Value = 2.86
$ws.Range("AA3").Value = 5.07
$ws.Range("AB3").Value = 4261.36
$ws.Range("AC3").Value = 691
$ws.Range("AD3").Value = 18.44
$ws.Range("AE3").Value = 23164
$ws.Range("AF3").Value = 0.55
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 1.57
$ws.Range("AI3").Value = 28.93
$ws.Range("AJ3").Value = 20000000

# Row 4
$ws.Range("D4").Value = 2149
$ws.Range("E4").Value = 96
$ws.Range("F4").Value = 96
$ws.Range("G4").Value = 159
$ws.Range("H4").Value = 132
$ws.Range("I4").Value = 132
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value = 4926
$ws.Range("L4").Value = 239
$ws.Range("M4").Value = 4686
$ws.Range("N4").Value = 4686
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = 100
$ws.Range("Q4").Value = 217
$ws.Range("R4").Value = -114
$ws.Range("S4").Value = -41
$ws.Range("T4").Value = 50
$ws.Range("U4").Value = 167
$ws.Range("V4").Value = 1
$ws.Range("W4").Value = 4.49
$ws.Range("X4").Value = 6.14
$ws.Range("Y4").Value = 2.84
$ws.Range("Z4").Value = 2.7
$ws.Range("AA4").Value = 5.11
$ws.Range("AB4").Value = 4341.28
$ws.Range("AC4").Value = 660
$ws.Range("AD4").Value = 18.4
$ws.Range("AE4").Value = 23431
$ws.Range("AF4").Value = 0.52
$ws.Range("AG4").Value = 200
$ws.Range("AH4").Value = 1.65
$ws.Range("AI4").Value = 30.29
$ws.Range("AJ4").Value = 20000000

# Row 5
$ws.Range("D5").Value = 2294
$ws.Range("E5").Value = 89
$ws.Range("F5").Value = 89
$ws.Range("G5").Value = 138
$ws.Range("H5").Value = 106
$ws.Range("I5").Value = 106
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value = 4976
$ws.Range("L5").Value = 268
$ws.Range("M5").Value = 4708
$ws.Range("N5").Value = 4708
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 100
$ws.Range("Q5").Value = 147
$ws.Range("R5").Value = -148
$ws.Range("S5").Value = -39
$ws.Range("T5").Value = 40
$ws.Range("U5").Value = 107
$ws.Range("V5").Value = 3
$ws.Range("W5").Value = 3.89
$ws.Range("X5").Value = 4.61
$ws.Range("Y5").Value = 2.25
$ws.Range("Z5").Value = 2.14
$ws.Range("AA5").Value = 5.68
$ws.Range("AB5").Value = 4398.31
$ws.Range("AC5").Value = 529
$ws.Range("AD5").Value = 24.67
$ws.Range("AE5").Value = 23542
$ws.Range("AF5").Value = 0.55
$ws.Range("AG5").Value = 250
$ws.Range("AH5").Value = 1.92
$ws.Range("AI5").Value = 47.26
$ws.Range("AJ5").Value = 20000000

# Row 6
$ws.Range("D6").Value = 2497
$ws.Range("E6").Value = 131
$ws.Range("F6").Value = 131
$ws.Range("G6").Value = 201
$ws.Range("H6").Value = 151
$ws.Range("I6").Value = 151
$ws.Range("K6").Value = 5080
$ws.Range("L6").Value = 276
$ws.Range("M6").Value = 4805
$ws.Range("N6").Value = 4805
$ws.Range("P6").Value = 100
$ws.Range("Q6").Value = 116
$ws.Range("R6").Value = -74
$ws.Range("S6").Value = -51
$ws.Range("T6").Value = 135
$ws.Range("U6").Value = -19
$ws.Range("V6").Value = 2
$ws.Range("W6").Value = 5.25
$ws.Range("X6").Value = 6.06
$ws.Range("Y6").Value = 3.18
$ws.Range("Z6").Value = 3.01
$ws.Range("AA6").Value = 5.73
$ws.Range("AB6").Value = 4510.58
$ws.Range("AC6").Value = 756
$ws.Range("AD6").Value = 15.66
$ws.Range("AE6").Value = 24024
$ws.Range("AF6").Value = 0.49
$ws.Range("AG6").Value = 250
$ws.Range("AH6").Value = 2.11
$ws.Range("AI6").Value = 33.05
$ws.Range("AJ6").Value = 20000000

# Row 7
$ws.Range("D7").Value = 2057
$ws.Range("E7").Value = 56
$ws.Range("G7").Value = 129
$ws.Range("H7").Value = 100
$ws.Range("I7").Value = 100
$ws.Range("K7").Value = 4892
$ws.Range("L7").Value = 232
$ws.Range("M7").Value = 4661
$ws.Range("N7").Value = 4661
$ws.Range("P7").Value = 100
$ws.Range("Q7").Value = 100
$ws.Range("R7").Value = -24
$ws.Range("S7").Value = -50
$ws.Range("T7").Value = 66
$ws.Range("U7").ClearContents()
$ws.Range("W7").Value = 2.72
$ws.Range("X7").Value = 4.86
$ws.Range("Y7").Value = 2.11
$ws.Range("Z7").Value = 2.01
$ws.Range("AA7").Value = 4.98
$ws.Range("AC7").Value = 500
$ws.Range("AD7").Value = 16.68
$ws.Range("AE7").Value = 23305
$ws.Range("AF7").Value = 0.36
$ws.Range("AG7").Value = 250
$ws.Range("AH7").Value = 3
$ws.Range("AI7").Value = 50

# Row 8
$ws.Range("D8").Value = 2077
$ws.Range("E8").Value = 62
$ws.Range("G8").Value = 143
$ws.Range("H8").Value = 111
$ws.Range("I8").Value = 111
$ws.Range("K8").Value = 4955
$ws.Range("L8").Value = 234
$ws.Range("M8").Value = 4722
$ws.Range("N8").Value = 4722
$ws.Range("P8").Value = 100
$ws.Range("Q8").Value = 144
$ws.Range("R8").Value = -78
$ws.Range("S8").Value = -50
$ws.Range("T8").Value = 70
$ws.Range("U8").ClearContents()
$ws.Range("W8").Value = 2.99
$ws.Range("X8").Value = 5.34
$ws.Range("Y8").Value = 2.37
$ws.Range("Z8").Value = 2.25
$ws.Range("AA8").Value = 4.96
$ws.Range("AC8").Value = 555
$ws.Range("AD8").Value = 15.03
$ws.Range("AE8").Value = 23610
$ws.Range("AF8").Value = 0.35
$ws.Range("AG8").Value = 250
$ws.Range("AH8").Value = 3
$ws.Range("AI8").Value = 45.05

# Row 9
$ws.Range("D9").Value = 2119
$ws.Range("E9").Value = 65
$ws.Range("G9").Value = 150
$ws.Range("H9").Value = 116
$ws.Range("I9").Value = 116
$ws.Range("K9").Value = 5025
$ws.Range("L9").Value = 238
$ws.Range("M9").Value = 4788
$ws.Range("N9").Value = 4788
$ws.Range("P9").Value = 100
$ws.Range("Q9").Value = 145
$ws.Range("R9").Value = -85
$ws.Range("S9").Value = -50
$ws.Range("T9").Value = 74
$ws.Range("U9").ClearContents()
$ws.Range("W9").Value = 3.07
$ws.Range("X9").Value = 5.47
$ws.Range("Y9").Value = 2.44
$ws.Range("Z9").Value = 2.33
$ws.Range("AA9").Value = 4.97
$ws.Range("AC9").Value = 580
$ws.Range("AD9").Value = 14.38
$ws.Range("AE9").Value = 23940
$ws.Range("AF9").Value = 0.35
$ws.Range("AG9").Value = 250
$ws.Range("AH9").Value = 3
$ws.Range("AI9").Value = 43.1
